$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data, preserving exact text formatting
# by forcing Text number format before assigning values that could otherwise
# be re-interpreted as numbers/dates by Excel (e.g. "1.00", "0.393", "7.40").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.775.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.386.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.388.46"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.393"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.966.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.392.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.874.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.526.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.69"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.417.69"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.96"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0776"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.556.34"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.30%  "
